$wb = $excel.ActiveWorkbook

$wsWebsiteRegistration = $wb.Worksheets.Item("websiteRegistarion")
$wsWebsiteRegistration.Range("A2").Value = "testweb30july22@gmail.com"

$wsAllreadyKyc = $wb.Worksheets.Item("allreadyKyc")
$wsAllreadyKyc.Range("A2").Value = "TEEPT2091J"
